$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 0.01372268907515325
$ws.Range("E2").Value = 0.3702421638625012
$ws.Range("I2").Value = 0.3991318617282426
$ws.Range("L2").Value = 0.5678023999999999
$ws.Range("M2").Value = 0.07761700000000001
$ws.Range("N2").Value = 12.63255121911827
$ws.Range("O2").Value = 3.433981830061676

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.04345357478734801
$ws.Range("E2").Value = 0.2300116167842811
$ws.Range("I2").Value = 0.5432101616050908
$ws.Range("L2").Value = 0.07828871198158027
$ws.Range("M2").Value = 0.05284941666666672
$ws.Range("N2").Value = 4.949474125454461
$ws.Range("O2").Value = 5.47229528149245

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.1185570572728702
$ws.Range("B2").Value = 0.03250564971082115
$ws.Range("E2").Value = 0.1631784079372048
$ws.Range("I2").Value = 0.4219666447996687
$ws.Range("M2").Value = 0.03633330321875103
$ws.Range("N2").Value = 7.447205196579223
$ws.Range("O2").Value = 0.8549932260662185

# --- Sheet "2040" ---
$ws = $wb.Worksheets.Item("2040")
$ws.Range("N2").Value = 0.5170046108778017
$ws.Range("O2").Value = 0

# --- Sheet "2045" ---
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 0.1328976602715719
$ws.Range("N2").Value = 2.090165151747748
$ws.Range("O2").Value = 4.911002349505852
